$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Pred Shifted = -4)
$ws.Range("C2").Value = -4.8372
$ws.Range("D2").Value = 0.659
$ws.Range("D2").Interior.Color = 10804650
$ws.Range("E2").Value = 2.1093
$ws.Range("E2").Interior.Color = 15136746

# Row 3 (Pred Shifted = -3)
$ws.Range("C3").Value = -0.3791
$ws.Range("D3").Value = 0.4876
$ws.Range("D3").Interior.Color = 6862424
$ws.Range("E3").Value = 1.1749
$ws.Range("E3").Interior.Color = 6599248

# Row 4 (Pred Shifted = -2)
$ws.Range("C4").Value = 0.5788
$ws.Range("D4").Value = 0.3598
$ws.Range("D4").Interior.Color = 4689446
$ws.Range("E4").Value = 0.8798
$ws.Range("E4").Interior.Color = 4623396

# Row 5 (Pred Shifted = -1)
$ws.Range("C5").Value = 0.9464
$ws.Range("D5").Value = 0.1605
$ws.Range("D5").Interior.Color = 1786880
$ws.Range("E5").Value = 0.3939
$ws.Range("E5").Interior.Color = 1786880

# Row 6 (Pred Shifted = 0)
$ws.Range("C6").Value = 0.6488
$ws.Range("D6").Value = 0.4591
$ws.Range("D6").Interior.Color = 6336072
$ws.Range("E6").Value = 1.1526
$ws.Range("E6").Interior.Color = 6467659

# Row 7 (Pred Shifted = 1)
$ws.Range("C7").Value = 0.3559
$ws.Range("D7").Value = 0.6243
$ws.Range("D7").Interior.Color = 9951132
$ws.Range("E7").Value = 1.5611
$ws.Range("E7").Interior.Color = 10082463

# Row 8 (Pred Shifted = 2)
$ws.Range("C8").Value = -0.0419
$ws.Range("D8").Value = 0.8072
$ws.Range("D8").Interior.Color = 14349279
$ws.Range("E8").Value = 1.9906
$ws.Range("E8").Interior.Color = 14283486

# Row 9 (Pred Shifted = 3)
$ws.Range("C9").Value = -0.3054
$ws.Range("D9").Value = 0.9215
$ws.Range("D9").Interior.Color = 16121079
$ws.Range("E9").Value = 2.2846
$ws.Range("E9").Interior.Color = 16121079
